$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 14 - "Crossed vs. Nested Design" - Content Placeholder 2
# Fix typo in the "nested" heuristic definition paragraph.
# ---------------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$sh14 = $s14.Shapes.Item(2)
$tr14 = $sh14.TextFrame.TextRange

$full = $tr14.Text
$needle = 'If factor '
$idx = $full.IndexOf($needle)
$tr14.Characters($idx + 1, $needle.Length).Text = 'A factor '

$full = $tr14.Text
$needle = 'in a factor A if every level of B occurs within exactly one level of A'
$idx = $full.IndexOf($needle)
$tr14.Characters($idx + 1, $needle.Length).Text = 'in a factor A if each level of B occurs in only one level of A'

# ---------------------------------------------------------------------------
# Slide 15 - "Crossed vs. Nested Design" (examples slide) - Content Placeholder 2
# ---------------------------------------------------------------------------
$s15 = $p.Slides.Item(15)
$sh15 = $s15.Shapes.Item(2)

# Nudge the placeholder up slightly (off.y 1693334 -> 1676400 EMU == 132.0pt)
$sh15.Top = 132.0

$tr15 = $sh15.TextFrame.TextRange

# --- paragraph 1: the quoted nested-factor definition -----------------
# "If factor B is " -> "A factor B is " and shrink to 18pt
$full = $tr15.Text
$needle = 'If factor B is '
$idx = $full.IndexOf($needle)
$tr15.Characters($idx + 1, $needle.Length).Text = 'A factor B is '
$tr15.Characters($idx + 1, 'A factor B is '.Length).Font.Size = 18

# "nested" -> shrink to 18pt
$full = $tr15.Text
$needle = 'nested'
$idx = $full.IndexOf($needle)
$tr15.Characters($idx + 1, $needle.Length).Font.Size = 18

# " in a factor A if each level of B occurs within only one level of " -> reworded, shrink to 18pt
$full = $tr15.Text
$needle = ' in a factor A if each level of B occurs within only one level of '
$idx = $full.IndexOf($needle)
$newText = ' in a factor A if each level of B occurs in only one level of '
$tr15.Characters($idx + 1, $needle.Length).Text = $newText
$tr15.Characters($idx + 1, $newText.Length).Font.Size = 18

# split trailing "A" off from the closing smart-quote and shrink just the "A" to 18pt
$full = $tr15.Text
$idx = $full.IndexOf('A"')
$tr15.Characters($idx + 1, 1).Font.Size = 18

# --- paragraph: "In this case the treatment is nested in gender." -----
$full = $tr15.Text
$oldLead = 'In this case the treatment is nested in gender'
$idx = $full.IndexOf($oldLead)
$tr15.Characters($idx + 1, $oldLead.Length).Text = 'In this case the treatment is nested in '

$full = $tr15.Text
$leadIdx = $full.IndexOf('In this case the treatment is nested in .')
$dotPos = $leadIdx + 'In this case the treatment is nested in '.Length
$tr15.Characters($dotPos + 1, 1).Text = 'gender'

$full = $tr15.Text
$genderIdx = $full.IndexOf('gender', $leadIdx)
$genderRng = $tr15.Characters($genderIdx + 1, 'gender'.Length)
$parenthetical = ' (in this case, gender is also nested in treatment, though it is more sensible the other way)'
$genderRng.InsertAfter($parenthetical) | Out-Null

$full = $tr15.Text
$tailIdx = $full.IndexOf('gender (in this case')
$tailStart = $tailIdx + 'gender'.Length
$tr15.Characters($tailStart + 1, $parenthetical.Length).Font.Size = 19

$full = $tr15.Text
$parenStart = $full.IndexOf('(in this case, gender is also nested')
$parenOnly = '(in this case, gender is also nested in treatment, though it is more sensible the other way)'
$tr15.Characters($parenStart + 1, $parenOnly.Length).Font.Size = 19
